$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.287.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.11%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.788.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.30%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

# Row 6: 'USDC' -> 'USDC'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

# Row 7: 'XRP' -> 'XRP'
$ws.Range("E7").Value = "  -2.52%  "

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3448"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "

# Row 9: 'OKB' -> 'OKB'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.10%  "

# Row 10: 'Polygon' -> 'Polygon'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.96%  "

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07405"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "

# Row 12: 'Solana' -> 'Solana'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.94%  "

# Row 13: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.09%  "

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.466"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "

# Row 15: 'Chainlink' -> 'Chainlink'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.369"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.02%  "

# Row 16: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.785.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

# Row 17: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001078"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.20%  "

# Row 18: 'TRON' -> 'TRON'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06674"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "

# Row 19: 'Litecoin' -> 'Litecoin'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "

# Row 20: 'Dai' -> 'Dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9989"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "

# Row 21: 'Avalanche' -> 'Avalanche'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "

# Row 22: 'Uniswap' -> 'Uniswap'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.462"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "

# Row 23: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.287.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.10%  "

# Row 24: 'Cosmos' -> 'Cosmos'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.97%  "

# Row 25: 'Toncoin' -> 'Toncoin'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.82%  "

# Row 26: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.445"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.45%  "

# Row 27: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.33%  "

# Row 28: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.431"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.82%  "

# Row 29: 'Monero' -> 'Monero'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.04%  "

# Row 30: 'BitcoinCash' -> 'WrappedliquidstakedEther2.0'
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.990.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "

# Row 31: 'WrappedliquidstakedEther2.0' -> 'BitcoinCash'
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.70%  "

# Row 32: 'Filecoin' -> 'HuobiToken'
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.020"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.33%  "

# Row 33: 'HuobiToken' -> 'Filecoin'
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.41%  "

# Row 34: 'Stellar' -> 'Stellar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08930"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.34%  "

# Row 35: 'Aptos' -> 'Aptos'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "

# Row 36: 'VeChain' -> 'VeChain'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02428"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.82%  "

# Row 37: 'TheSandbox' -> 'TheSandbox'
$ws.Range("E37").Value = "  +1.48%  "

# Row 38: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.39%  "

# Row 39: 'Hedera' -> 'Hedera'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06402"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "

# Row 40: 'Algorand' -> 'Algorand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "

# Row 41: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.248"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "

# Row 42: 'WEMIXTOKEN' -> 'WEMIXTOKEN'
$ws.Range("E42").Value = "  -6.69%  "

# Row 43: 'FraxShare' -> 'FraxShare'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.83%  "

# Row 44: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "

# Row 45: 'Frax' -> 'Frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46: 'Decentraland' -> 'Decentraland'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6327"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "

# Row 47: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.879"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.64%  "

# Row 48: 'Quant' -> 'Quant'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.04%  "

# Row 49: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.086"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.23%  "

# Row 50: 'Cronos' -> 'Cronos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07497"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.61%  "

# Row 51: 'EOS' -> 'EOS'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.212"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.72%  "
